$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-9 (id, question, op1, op2, op3, op4, ans)
$data = @(
    @(1, "Solve for x: 15x = 75", 3, 5, 6, 10, 5),
    @(2, "Solve for x: 5x - 2 = 13", 1, 2, 3, 4, 3),
    @(3, "Solve for x: 4x + 6 = 18", 1, 2, 3, 4, 3),
    @(4, "Solve for x: 6x = 42", 6, 7, 8, 12, 7),
    @(5, "Solve for x: 3x + 5 = 14", 2, 3, 4, 5, 3),
    @(6, "Solve for x: 6x - 9 = 15", 2, 3, 4, 5, 4),
    @(7, "Solve for x: 4x + 5 = 2x + 11", 1, 2, 3, 4, 3),
    @(8, "Solve for x: 6x + 3 = 4x + 9", 2, 3, 4, 5, 3)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $ws.Cells.Item($row, 6).Value = $entry[5]
    $ws.Cells.Item($row, 7).Value = $entry[6]
    $row++
}
